# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.108.45"
$ws.Range("E2").Value = "  -0.54%  "

$ws.Range("D3").Value = "3.133.75"
$ws.Range("E3").Value = "  -1.39%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.42"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.97"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.11%  "

$ws.Range("E7").Value = "  -0.24%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.571"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -6.04%  "

$ws.Range("D9").Value = "3.144.06"
$ws.Range("E9").Value = "  -1.11%  "

$ws.Range("E10").Value = "  -3.79%  "

$ws.Range("E11").Value = "  -3.43%  "

$ws.Range("E12").Value = "  -0.70%  "

$ws.Range("D13").Value = "3.678.02"
$ws.Range("E13").Value = "  -1.31%  "

$ws.Range("E14").Value = "  -0.69%  "

$ws.Range("D15").Value = "64.142.82"
$ws.Range("E15").Value = "  -0.57%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "24.99"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.51%  "

$ws.Range("D17").Value = "3.144.44"
$ws.Range("E17").Value = "  -1.03%  "

$ws.Range("E18").Value = "  -3.80%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "398.92"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.80%  "

$ws.Range("E20").Value = "  -3.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.47"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.08"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.26%  "

$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("E24").Value = "  -2.81%  "

$ws.Range("E25").Value = "  -1.28%  "

$ws.Range("E26").Value = "  -5.82%  "

$ws.Range("E27").Value = "  -5.45%  "

$ws.Range("E28").Value = "  -1.64%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.995"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.36%  "

$ws.Range("E30").Value = "  +0.02%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.79"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.24%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.07"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.92%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "159.99"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.40%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.24"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.63%  "

$ws.Range("E35").Value = "  -4.80%  "

$ws.Range("E36").Value = "  -3.82%  "

$ws.Range("E37").Value = "  -2.99%  "

$ws.Range("D38").Value = "2.642.59"
$ws.Range("E38").Value = "  -3.41%  "

$ws.Range("E39").Value = "  -2.95%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.45"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.55%  "

$ws.Range("E41").Value = "  -3.23%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.21"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.47%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.687"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.26%  "

$ws.Range("E44").Value = "  -2.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.38"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.40%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0253"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.57%  "

$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "284.33"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.35%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.94"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.18%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.996"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.36%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0970"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.48"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.37%  "
